$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the formatting of the other
# header cells (B1:G1) by copying G1's format onto H1, then setting
# the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H8 with the value 1 (plain numeric cells, no special style)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
